$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 300
$ws.Range("I38").Value = 100
$ws.Range("J38").Value = 500
$ws.Range("K38").Value = 300
$ws.Range("L38").Value = 1500
$ws.Range("M38").Value = 72
$ws.Range("N38").Value = -2244

$ws.Range("H58").Value = 122.1
$ws.Range("J58").Value = 166.75
$ws.Range("L58").Value = 500.25
$ws.Range("N58").Value = -800.25

$ws.Range("I106").Value = 2668
$ws.Range("K106").Value = 2668
$ws.Range("M106").Value = -2037

$ws.Range("H107").Value = 1027.3334
$ws.Range("I107").Value = 417.14285
$ws.Range("K107").Value = 417.14285
$ws.Range("M107").Value = 1502.85715

$ws.Range("H132").Value = 4353.353
$ws.Range("I132").Value = 3642.9285
$ws.Range("J132").Value = 7668.6665
$ws.Range("K132").Value = 10928.7855
$ws.Range("L132").Value = 23005.9995
$ws.Range("M132").Value = -8398.7855
$ws.Range("N132").Value = -28065.9995

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1050
$ws.Range("I2").Value = 600
$ws.Range("J2").Value = 1500
$ws.Range("K2").Value = 600
$ws.Range("L2").Value = 1500
$ws.Range("M2").Value = -487
$ws.Range("N2").Value = -1726

$ws.Range("H45").Value = 1216.6666
$ws.Range("I45").Value = 1216.6666
$ws.Range("K45").Value = 1216.6666
$ws.Range("M45").Value = -839.6666

$ws.Range("H74").Value = 3274.7
$ws.Range("I74").Value = 1249.5714
$ws.Range("J74").Value = 8000
$ws.Range("K74").Value = 1249.5714
$ws.Range("L74").Value = 8000
$ws.Range("M74").Value = -375.5714
$ws.Range("N74").Value = -9748

$ws.Range("H77").Value = 3274.7
$ws.Range("I77").Value = 1249.5714
$ws.Range("J77").Value = 8000
$ws.Range("K77").Value = 6247.857
$ws.Range("L77").Value = 40000
$ws.Range("M77").Value = -1879.857
$ws.Range("N77").Value = -48736

$ws.Range("H103").Value = 20000
$ws.Range("J103").Value = 20000
$ws.Range("L103").Value = 20000
$ws.Range("N103").Value = -22344

$ws.Range("H110").Value = 920
$ws.Range("J110").Value = 750
$ws.Range("L110").Value = 750
$ws.Range("N110").Value = -4840

$ws.Range("H116").Value = 1050
$ws.Range("I116").Value = 600
$ws.Range("J116").Value = 1500
$ws.Range("K116").Value = 600
$ws.Range("L116").Value = 1500
$ws.Range("M116").Value = 1694
$ws.Range("N116").Value = -6088

$ws.Range("H122").Value = 1836.1333
$ws.Range("I122").Value = 1770.1666
$ws.Range("K122").Value = 5310.4998
$ws.Range("M122").Value = -2860.4998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1050
$ws.Range("I3").Value = 600
$ws.Range("J3").Value = 1500
$ws.Range("K3").Value = 600
$ws.Range("L3").Value = 1500
$ws.Range("M3").Value = -486
$ws.Range("N3").Value = -1728

$ws.Range("H86").Value = 1449.5
$ws.Range("J86").Value = 1400
$ws.Range("L86").Value = 1400
$ws.Range("N86").Value = -3646

$ws.Range("H89").Value = 1449.5
$ws.Range("J89").Value = 1400
$ws.Range("L89").Value = 7000
$ws.Range("N89").Value = -18232

$ws.Range("H105").Value = 3857.25
$ws.Range("I105").Value = 3189.75
$ws.Range("J105").Value = 4524.75
$ws.Range("K105").Value = 3189.75
$ws.Range("L105").Value = 4524.75
$ws.Range("M105").Value = -1442.75
$ws.Range("N105").Value = -8018.75

$ws.Range("H107").Value = 1338.2
$ws.Range("I107").Value = 1034.5714
$ws.Range("K107").Value = 1034.5714
$ws.Range("M107").Value = 885.4286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2532.7778
$ws.Range("I122").Value = 1849.6666
$ws.Range("K122").Value = 5548.9998
$ws.Range("M122").Value = -3098.9998

$ws.Range("H132").Value = 3825
$ws.Range("I132").Value = 3952.75
$ws.Range("K132").Value = 11858.25
$ws.Range("M132").Value = -9328.25

$ws.Range("H134").Value = 1009
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 900
$ws.Range("I26").Value = 900
$ws.Range("J26").Value = 900
$ws.Range("K26").Value = 2700
$ws.Range("L26").Value = 2700
$ws.Range("M26").Value = -2412
$ws.Range("N26").Value = -3276

$ws.Range("H34").Value = 232.66667
$ws.Range("J34").Value = 251.5
$ws.Range("L34").Value = 754.5
$ws.Range("N34").Value = -922.5

$ws.Range("H103").Value = 1672.6
$ws.Range("I103").Value = 1300.1666
$ws.Range("J103").Value = 2231.25
$ws.Range("K103").Value = 3900.4998
$ws.Range("L103").Value = 6693.75
$ws.Range("M103").Value = -3021.4998
$ws.Range("N103").Value = -8451.75

$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws.Range("H131").Value = 2385.5264
$ws.Range("I131").Value = 1480
$ws.Range("J131").Value = 2492.0588
$ws.Range("K131").Value = 4440
$ws.Range("L131").Value = 7476.176399999999
$ws.Range("M131").Value = 600
$ws.Range("N131").Value = -17556.1764

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1842.8572
$ws.Range("I80").Value = 1300.3334
$ws.Range("J80").Value = 2249.75
$ws.Range("K80").Value = 1300.3334
$ws.Range("L80").Value = 2249.75
$ws.Range("M80").Value = -302.3334
$ws.Range("N80").Value = -4245.75

$ws.Range("H83").Value = 1842.8572
$ws.Range("I83").Value = 1300.3334
$ws.Range("J83").Value = 2249.75
$ws.Range("K83").Value = 6501.666999999999
$ws.Range("L83").Value = 11248.75
$ws.Range("M83").Value = -1509.666999999999
$ws.Range("N83").Value = -21232.75

$ws.Range("H97").Value = 916
$ws.Range("I97").Value = 916
$ws.Range("K97").Value = 916
$ws.Range("M97").Value = -420

$ws.Range("H102").Value = 1969.9615
$ws.Range("I102").Value = 2052.9048
$ws.Range("J102").Value = 1621.6
$ws.Range("K102").Value = 2052.9048
$ws.Range("L102").Value = 1621.6
$ws.Range("M102").Value = -430.9047999999998
$ws.Range("N102").Value = -4865.6

$ws.Range("H122").Value = 1071.9
$ws.Range("I122").Value = 1071.9
$ws.Range("K122").Value = 3215.7
$ws.Range("M122").Value = -765.7000000000003

$ws.Range("H132").Value = 4874.4614
$ws.Range("I132").Value = 4447.3335
$ws.Range("K132").Value = 13342.0005
$ws.Range("M132").Value = -10812.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1578.3684
$ws.Range("I22").Value = 1606.5714
$ws.Range("J22").Value = 1499.4
$ws.Range("K22").Value = 1606.5714
$ws.Range("L22").Value = 1499.4
$ws.Range("M22").Value = -1311.5714
$ws.Range("N22").Value = -2089.4

$ws.Range("H27").Value = 1578.3684
$ws.Range("I27").Value = 1606.5714
$ws.Range("J27").Value = 1499.4
$ws.Range("K27").Value = 1606.5714
$ws.Range("L27").Value = 1499.4
$ws.Range("M27").Value = -1499.5714
$ws.Range("N27").Value = -1713.4

$ws.Range("H46").Value = 4164.436
$ws.Range("I46").Value = 2189.9
$ws.Range("K46").Value = 2189.9
$ws.Range("M46").Value = -2001.9

$ws.Range("H55").Value = 2562.875
$ws.Range("I55").Value = 2071.8572
$ws.Range("K55").Value = 2071.8572
$ws.Range("M55").Value = -1898.8572

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H100").Value = 11250
$ws.Range("I100").Value = 20000
$ws.Range("K100").Value = 20000
$ws.Range("M100").Value = -19459

$ws.Range("H101").Value = 11000.333
$ws.Range("J101").Value = 11000.333
$ws.Range("L101").Value = 11000.333
$ws.Range("N101").Value = -17490.333

$ws.Range("H136").Value = 4787.6665
$ws.Range("I136").Value = 4545.2
$ws.Range("K136").Value = 13635.6
$ws.Range("M136").Value = -11085.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 18000
$ws.Range("I28").Value = 18000
$ws.Range("K28").Value = 18000
$ws.Range("M28").Value = -17652

$ws.Range("H92").Value = 20000
$ws.Range("J92").Value = 20000
$ws.Range("L92").Value = 20000
$ws.Range("N92").Value = -24992

$ws.Range("H122").Value = 1864.1666
$ws.Range("I122").Value = 1864.1666
$ws.Range("K122").Value = 5592.4998
$ws.Range("M122").Value = -3142.4998

$ws.Range("H126").Value = 4612.25
$ws.Range("J126").Value = 3966.6667
$ws.Range("L126").Value = 11900.0001
$ws.Range("N126").Value = -16840.0001
